$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Combined")
$ws.Activate()

$ws.Range("B1").Value = "SPF CPI"
$ws.Range("C1").Value = "SPF PCE"
$ws.Range("D1").Value = "SCE"
$ws.Range("A2").Value = "Test 1: Bias"
$ws.Range("A3").Value = "Constant"
$ws.Range("B3").Formula = "=SPFInd!B20"
$ws.Range("C3").Formula = "=SPFInd!G5"
$ws.Range("D3").Formula = "=SCEInd!B71"
$ws.Range("B4").Formula = "=SPFInd!B21"
$ws.Range("C4").Formula = "=SPFInd!G6"
$ws.Range("D4").Formula = "=SCEInd!B72"
$ws.Range("A5").Value = "N"
$ws.Range("B5").Formula = "=SPFInd!B23"
$ws.Range("C5").Formula = "=SPFInd!G23"
$ws.Range("D5").Formula = "=SCEInd!B74"
$ws.Range("A8").Value = "Test2: FE Depends on past information"
$ws.Range("A9").Value = "Forecast 1-yr before"
$ws.Range("B9").Formula = "=SPFInd!C5"
$ws.Range("C9").Formula = "=SPFInd!G5"
$ws.Range("D9").Value = "NA"
$ws.Range("B10").Formula = "=SPFInd!C6"
$ws.Range("C10").Formula = "=SPFInd!G6"
$ws.Range("D10").Value = "NA"
$ws.Range("A11").Value = "Constant"
$ws.Range("B11").Formula = "=SPFInd!C20"
$ws.Range("C11").Formula = "=SPFInd!G20"
$ws.Range("D11").Value = "NA"
$ws.Range("B12").Formula = "=SPFInd!C21"
$ws.Range("C12").Formula = "=SPFInd!G21"
$ws.Range("D12").Value = "NA"
$ws.Range("A13").Value = "N"
$ws.Range("B13").Formula = "=SPFInd!C23"
$ws.Range("C13").Formula = "=SPFInd!G23"
$ws.Range("D13").Value = "NA"
$ws.Range("A14").Value = "R^2"
$ws.Range("B14").Formula = "=SPFInd!C24"
$ws.Range("C14").Formula = "=SPFInd!G24"
$ws.Range("D14").Value = "NA"
$ws.Range("A16").Value = "Test3: FE of non-overllaping forecast horizons are serially correlated "
$ws.Range("A17").Value = "Forecast Error 1-year before"
$ws.Range("B17").Formula = "=SPFInd!D8"
$ws.Range("C17").Formula = "=SPFInd!H8"
$ws.Range("D17").Value = "NA"
$ws.Range("B18").Formula = "=SPFInd!D9"
$ws.Range("C18").Formula = "=SPFInd!H9"
$ws.Range("D18").Value = "NA"
$ws.Range("A19").Value = "Constant"
$ws.Range("B19").Formula = "=SPFInd!D20"
$ws.Range("C19").Formula = "=SPFInd!H20"
$ws.Range("D19").Value = "NA"
$ws.Range("B20").Formula = "=SPFInd!D21"
$ws.Range("C20").Formula = "=SPFInd!H21"
$ws.Range("D20").Value = "NA"
$ws.Range("D21").Value = "NA"
$ws.Range("A22").Value = "N"
$ws.Range("B22").Formula = "=SPFInd!D23"
$ws.Range("C22").Formula = "=SPFInd!H23"
$ws.Range("D22").Value = "NA"
$ws.Range("A23").Value = "R^2"
$ws.Range("B23").Formula = "=SPFInd!D24"
$ws.Range("C23").Formula = "=SPFInd!H24"
$ws.Range("D23").Value = "NA"
$ws.Range("A26").Value = "Test4: Overlapping FE are serially correlated "
$ws.Range("A28").Value = "Forecast Error 1-q before"
$ws.Range("B28").Formula = "=SPFInd!E11"
$ws.Range("C28").Formula = "=SPFInd!I11"
$ws.Range("D28").Formula = "=SCEInd!D14"
$ws.Range("B29").Formula = "=SPFInd!E12"
$ws.Range("C29").Formula = "=SPFInd!I12"
$ws.Range("D29").Formula = "=SCEInd!D15"
$ws.Range("A30").Value = "Forecast Error 2-q before"
$ws.Range("B30").Formula = "=SPFInd!E14"
$ws.Range("C30").Formula = "=SPFInd!I14"
$ws.Range("D30").Formula = "=SCEInd!E26"
$ws.Range("B31").Formula = "=SPFInd!E15"
$ws.Range("C31").Formula = "=SPFInd!I15"
$ws.Range("D31").Formula = "=SCEInd!E27"
$ws.Range("A32").Value = "Forecast Error 3-q before"
$ws.Range("B32").Formula = "=SPFInd!E17"
$ws.Range("C32").Formula = "=SPFInd!I17"
$ws.Range("D32").Formula = "=SCEInd!H59"
$ws.Range("B33").Formula = "=SPFInd!E18"
$ws.Range("C33").Formula = "=SPFInd!I18"
$ws.Range("D33").Formula = "=SCEInd!H60"
$ws.Range("A35").Value = "Constant"
$ws.Range("B35").Formula = "=SPFInd!E20"
$ws.Range("C35").Formula = "=SPFInd!I20"
$ws.Range("D35").Formula = "=SCEInd!H71"
$ws.Range("B36").Formula = "=SPFInd!E21"
$ws.Range("C36").Formula = "=SPFInd!I21"
$ws.Range("D36").Formula = "=SCEInd!H72"
$ws.Range("A38").Value = "N"
$ws.Range("B38").Formula = "=SPFInd!E23"
$ws.Range("C38").Formula = "=SPFInd!I23"
$ws.Range("D38").Formula = "=SCEInd!H74"
$ws.Range("A39").Value = "R^2"
$ws.Range("B39").Formula = "=SPFInd!E24"
$ws.Range("C39").Formula = "=SPFInd!I24"
$ws.Range("D39").Formula = "=SCEInd!H75"

# Wrap text + row height 48 for style-1 rows (17,28,29,30,31,32)
$ws.Rows.Item(17).RowHeight = 48
$ws.Rows.Item(28).RowHeight = 48
$ws.Rows.Item(29).RowHeight = 48
$ws.Rows.Item(30).RowHeight = 48
$ws.Rows.Item(31).RowHeight = 48
$ws.Rows.Item(32).RowHeight = 48
$ws.Range("A17").WrapText = $true
$ws.Range("A28").WrapText = $true
$ws.Range("A29").WrapText = $true
$ws.Range("A30").WrapText = $true
$ws.Range("A31").WrapText = $true
$ws.Range("A32").WrapText = $true

# Page setup: portrait orientation
$ws.PageSetup.Orientation = 1

# Selection on SCEInd sheet (view state)
$wsSCE = $wb.Worksheets.Item("SCEInd")
$wsSCE.Activate()
$wsSCE.Range("E15").Select()

# Return to Combined as the active sheet with final selection
$ws.Activate()
$ws.Range("F35").Select()
